# Title screen + sound settings
# - HUD and UI row: now meets standard, description changed to the
#   player-info/menu-system text
# - Sound row: now present (1), with a "Background music" description
# - The "Basic Game rubric" sheet becomes the active/selected sheet again

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Basic Game rubric")

$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "Show Player info + Menu system"

$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "Background music"

$ws.Select()
